# Add "setup00" / "Software installation" row to the syllabus,
# inserted above the existing first schedule entry (row 2), pushing
# everything else down by one row. Excel auto-adjusts the relative
# date formulas (B col) and shared-formula ranges when the row shifts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2 (shifts rows 2..21 down to 3..22,
# and relative formulas shift/re-anchor automatically).
$ws.Rows.Item(2).Insert()

# Populate the new row's topic / flag / title. Leave the date (column B)
# empty for this row -- there is no scheduled date for it.
$ws.Range("A2").Value = "setup00"
$ws.Range("C2").Value = $true
$ws.Range("D2").Value = "Software installation"
$ws.Range("B2").Clear()

# Match the author's final selection state.
$ws.Range("A2").Select()
